$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Make G10 match the formatting of the rest of the "durumu" column (G4:G9)
# before writing the formula into it (mirrors what a user gets after
# selecting G4:G10 and filling the IF() formula down the whole range).
$ws.Range("G9").Copy()
$ws.Range("G10").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

# Durumu column: IF(not<45,"KALDI","GEÇTİ") for every student row
$ws.Range("G4").Formula = '=IF(F4<45,"KALDI","GEÇTİ")'
$ws.Range("G5:G10").Formula = '=IF(F5<45,"KALDI","GEÇTİ")'

# Fill in the student identification block (Numara / Ad Soyad / Bölüm)
$ws.Range("L8").Value = 20215070019
$ws.Range("L9").Value = "KÜBRA ÇABUK"
$ws.Range("L10").Value = "YBS"

# Leave the same cell selected as in the edited workbook
$ws.Range("L10:N10").Select() | Out-Null
